$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 90.75
$ws.Range("I33").Value = 94.5
$ws.Range("J33").Value = 57
$ws.Range("K33").Value = 94.5
$ws.Range("L33").Value = 57
$ws.Range("M33").Value = 134.5
$ws.Range("N33").Value = -515
$ws.Range("H64").Value = 3659.6667
$ws.Range("J64").Value = 3990
$ws.Range("L64").Value = 3990
$ws.Range("N64").Value = -4486
$ws.Range("H67").Value = 3659.6667
$ws.Range("J67").Value = 3990
$ws.Range("L67").Value = 3990
$ws.Range("N67").Value = -5706
$ws.Range("H76").Value = 7872.778
$ws.Range("J76").Value = 7835.5713
$ws.Range("L76").Value = 7835.5713
$ws.Range("N76").Value = -8465.5713
$ws.Range("H79").Value = 7872.778
$ws.Range("J79").Value = 7835.5713
$ws.Range("L79").Value = 7835.5713
$ws.Range("N79").Value = -10019.5713
$ws.Range("H96").Value = 1481.1
$ws.Range("I96").Value = 1481.1
$ws.Range("K96").Value = 4443.299999999999
$ws.Range("M96").Value = -3070.299999999999
$ws.Range("H137").Value = 14432.762
$ws.Range("I137").Value = 9209.799999999999
$ws.Range("J137").Value = 19180.908
$ws.Range("K137").Value = 27629.4
$ws.Range("L137").Value = 57542.724
$ws.Range("M137").Value = -25079.4
$ws.Range("N137").Value = -62642.724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 304.75
$ws.Range("I4").Value = 273
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 273
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -157
$ws.Range("N4").Value = -632
$ws.Range("H36").Value = 7173.6665
$ws.Range("I36").Value = 7173.6665
$ws.Range("K36").Value = 7173.6665
$ws.Range("M36").Value = -6827.6665
$ws.Range("H37").Value = 13999.5
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 13999.5
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 13999.5
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -14545.5
$ws.Range("H63").Value = 6721.8125
$ws.Range("I63").Value = 1849.6666
$ws.Range("J63").Value = 7846.154
$ws.Range("K63").Value = 1849.6666
$ws.Range("L63").Value = 7846.154
$ws.Range("M63").Value = -1163.6666
$ws.Range("N63").Value = -9218.154
$ws.Range("H64").Value = 65000
$ws.Range("J64").Value = 65000
$ws.Range("L64").Value = 65000
$ws.Range("N64").Value = -65496
$ws.Range("H66").Value = 6721.8125
$ws.Range("I66").Value = 1849.6666
$ws.Range("J66").Value = 7846.154
$ws.Range("K66").Value = 9248.333000000001
$ws.Range("L66").Value = 39230.77
$ws.Range("M66").Value = -5816.333000000001
$ws.Range("N66").Value = -46094.77
$ws.Range("H67").Value = 65000
$ws.Range("J67").Value = 65000
$ws.Range("L67").Value = 65000
$ws.Range("N67").Value = -66716
$ws.Range("H80").Value = 94999.25
$ws.Range("J80").Value = 94999.25
$ws.Range("L80").Value = 94999.25
$ws.Range("N80").Value = -96995.25
$ws.Range("H83").Value = 94999.25
$ws.Range("J83").Value = 94999.25
$ws.Range("L83").Value = 284997.75
$ws.Range("N83").Value = -294981.75
$ws.Range("H97").Value = 524.8889
$ws.Range("I97").Value = 549.875
$ws.Range("J97").Value = 325
$ws.Range("K97").Value = 549.875
$ws.Range("L97").Value = 325
$ws.Range("M97").Value = -53.875
$ws.Range("N97").Value = -1317
$ws.Range("H122").Value = 4647.9
$ws.Range("I122").Value = 3993.3333
$ws.Range("K122").Value = 11979.9999
$ws.Range("M122").Value = -9529.999899999999
$ws.Range("H132").Value = 3743.8333
$ws.Range("I132").Value = 2489.8
$ws.Range("K132").Value = 7469.400000000001
$ws.Range("M132").Value = -4939.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 86657.836
$ws.Range("I22").Value = 204153.8
$ws.Range("J22").Value = 2732.1428
$ws.Range("K22").Value = 204153.8
$ws.Range("L22").Value = 2732.1428
$ws.Range("M22").Value = -203803.8
$ws.Range("N22").Value = -3432.1428
$ws.Range("H31").Value = 9438.166999999999
$ws.Range("I31").Value = 13448.375
$ws.Range("J31").Value = 6230
$ws.Range("K31").Value = 13448.375
$ws.Range("L31").Value = 6230
$ws.Range("M31").Value = -13153.375
$ws.Range("N31").Value = -6820
$ws.Range("H34").Value = 9438.166999999999
$ws.Range("I34").Value = 13448.375
$ws.Range("J34").Value = 6230
$ws.Range("K34").Value = 13448.375
$ws.Range("L34").Value = 6230
$ws.Range("M34").Value = -13246.375
$ws.Range("N34").Value = -6634
$ws.Range("H80").Value = 59750
$ws.Range("J80").Value = 59750
$ws.Range("L80").Value = 59750
$ws.Range("N80").Value = -61996
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H83").Value = 59750
$ws.Range("J83").Value = 59750
$ws.Range("L83").Value = 179250
$ws.Range("N83").Value = -190482
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H134").Value = 1458.6666
$ws.Range("I134").Value = 1029.4615
$ws.Range("K134").Value = 3088.3845
$ws.Range("M134").Value = -553.3844999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 76.333336
$ws.Range("I60").Value = 76.333336
$ws.Range("K60").Value = 229.000008
$ws.Range("M60").Value = 21.99999199999999
$ws.Range("H80").Value = 2998.5
$ws.Range("I80").Value = 2998
$ws.Range("K80").Value = 8994
$ws.Range("M80").Value = -8058
$ws.Range("H83").Value = 2998.5
$ws.Range("I83").Value = 2998
$ws.Range("K83").Value = 26982
$ws.Range("M83").Value = -22302

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H74").Value = 42500
$ws.Range("I74").Value = 30000
$ws.Range("K74").Value = 30000
$ws.Range("M74").Value = -29064
$ws.Range("H77").Value = 42500
$ws.Range("I77").Value = 30000
$ws.Range("K77").Value = 90000
$ws.Range("M77").Value = -85320
$ws.Range("H80").Value = 3100
$ws.Range("I80").Value = 3071.4285
$ws.Range("J80").Value = 3166.6667
$ws.Range("K80").Value = 3071.4285
$ws.Range("L80").Value = 3166.6667
$ws.Range("M80").Value = -2073.4285
$ws.Range("N80").Value = -5162.6667
$ws.Range("H83").Value = 3100
$ws.Range("I83").Value = 3071.4285
$ws.Range("J83").Value = 3166.6667
$ws.Range("K83").Value = 15357.1425
$ws.Range("L83").Value = 15833.3335
$ws.Range("M83").Value = -10365.1425
$ws.Range("N83").Value = -25817.3335
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4218.8335
$ws.Range("I7").Value = 4218.8335
$ws.Range("K7").Value = 4218.8335
$ws.Range("M7").Value = -4106.8335
$ws.Range("H34").Value = 24874.5
$ws.Range("I34").Value = 29999.334
$ws.Range("K34").Value = 29999.334
$ws.Range("M34").Value = -29827.334
$ws.Range("H46").Value = 3448
$ws.Range("I46").Value = 1125
$ws.Range("J46").Value = 4996.6665
$ws.Range("K46").Value = 1125
$ws.Range("L46").Value = 4996.6665
$ws.Range("M46").Value = -937
$ws.Range("N46").Value = -5372.6665
$ws.Range("H68").Value = 4449.5
$ws.Range("I68").Value = 4131.25
$ws.Range("K68").Value = 4131.25
$ws.Range("M68").Value = -3382.25
$ws.Range("H71").Value = 4449.5
$ws.Range("I71").Value = 4131.25
$ws.Range("K71").Value = 20656.25
$ws.Range("M71").Value = -16912.25
$ws.Range("H82").Value = 201059.2
$ws.Range("I82").Value = 1515.6666
$ws.Range("J82").Value = 500374.5
$ws.Range("K82").Value = 1515.6666
$ws.Range("L82").Value = 500374.5
$ws.Range("M82").Value = -1154.6666
$ws.Range("N82").Value = -501096.5
$ws.Range("H85").Value = 201059.2
$ws.Range("I85").Value = 1515.6666
$ws.Range("J85").Value = 500374.5
$ws.Range("K85").Value = 1515.6666
$ws.Range("L85").Value = 500374.5
$ws.Range("M85").Value = -267.6666
$ws.Range("N85").Value = -502870.5
$ws.Range("H100").Value = 7666
$ws.Range("I100").Value = 5571.2856
$ws.Range("J100").Value = 14997.5
$ws.Range("K100").Value = 5571.2856
$ws.Range("L100").Value = 14997.5
$ws.Range("M100").Value = -5030.2856
$ws.Range("N100").Value = -16079.5
$ws.Range("H126").Value = 4218.8335
$ws.Range("I126").Value = 4218.8335
$ws.Range("K126").Value = 12656.5005
$ws.Range("M126").Value = -10186.5005
$ws.Range("H140").Value = 29100
$ws.Range("J140").Value = 36520
$ws.Range("L140").Value = 36520
$ws.Range("N140").Value = -46880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3800.8333
$ws.Range("I96").Value = 3561
$ws.Range("K96").Value = 3561
$ws.Range("M96").Value = -2188
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 2414.3333
$ws.Range("I107").Value = 1956.8
$ws.Range("J107").Value = 3558.1667
$ws.Range("K107").Value = 5870.4
$ws.Range("L107").Value = 10674.5001
$ws.Range("M107").Value = -3950.4
$ws.Range("N107").Value = -14514.5001

